$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All existing data rows (2-329) had their "Förändrad" (C column) date
# bumped from 2023-09-13 (45182) to 2023-09-15 (45184).
$ws.Range("C2:C329").Value = 45184

# Row 329 picks up an explicit row height (ht="15" customHeight="1").
$ws.Rows.Item(329).RowHeight = 15

# Copy the formatting (date number format / wrap-text style) of row 329
# down into the brand-new row 330 before filling in its values.
$ws.Range("B329").Copy($ws.Range("B330"))
$ws.Range("C329").Copy($ws.Range("C330"))
$ws.Range("R329").Copy($ws.Range("R330"))

# New record appended at the bottom of the table.
$ws.Range("A330").Value = "A 43146-2023"
$ws.Range("B330").Value = 45183
$ws.Range("C330").Value = 45184
$ws.Range("D330").Value = "BLEKINGE LÄN"
$ws.Range("E330").Value = "KARLSHAMN"
$ws.Range("G330").Value = 1.4
$ws.Range("H330").Value = 0
$ws.Range("I330").Value = 0
$ws.Range("J330").Value = 0
$ws.Range("K330").Value = 0
$ws.Range("L330").Value = 0
$ws.Range("M330").Value = 0
$ws.Range("N330").Value = 0
$ws.Range("O330").Value = 0
$ws.Range("P330").Value = 0
$ws.Range("Q330").Value = 0
$ws.Range("R330").Value = ""
